# Insert two new rows before the existing row 13 (pushing the old
# "Cryptography / IoT / MiniProject" rows down to 15-17) and populate the
# two new rows with the "Value added course" and "Spoken Tutorial" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("13:14").Insert()

# New row 13: Value added course - Entrepreneurship Development
$ws.Range("A13").Value = "VI CSE"
$ws.Range("B13").Value = "ED2VA1"
$ws.Range("C13").Value = "Value added course - Entrepreneurship Development"
$ws.Range("D13").Value = "PAC"
$ws.Range("E13").Value = "NPP"
$ws.Range("F13").Value = "SGR"
$ws.Range("G13").Value = "RSA"
# H13 is an explicit empty-text cell (matches the other rows' blank "E" column).
# A plain "" assignment clears the cell instead of leaving an empty string, so
# use the leading-apostrophe text marker and then strip the quote-prefix style
# it leaves behind.
$ws.Range("H13").Value = "'"
$ws.Range("H13").ClearFormats()
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = "CSE"
$ws.Range("K13").Value = "NO"
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 0

# New row 14: IIT Spoken Tutorial class
$ws.Range("A14").Value = "VI CSE"
$ws.Range("B14").Value = "IT3412"
$ws.Range("C14").Value = "IIT Spoken Tutorial class"
$ws.Range("D14").Value = "LA1"
$ws.Range("E14").Value = "LA2"
$ws.Range("F14").Value = "LA3"
$ws.Range("G14").Value = "LA4"
$ws.Range("H14").Value = "'"
$ws.Range("H14").ClearFormats()
$ws.Range("I14").Value = 4
$ws.Range("J14").Value = "IT"
$ws.Range("K14").Value = "NO"
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 1
